$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MEG_ANALYSIS_MASTERFILE")

# Rows where the ANATOMY block (columns C:H) goes from "only nifti flag known"
# to "nifti collected (=1), recon/bem/coreg/src stages now tracked as 0" --
# i.e. column C's highlight is cleared and column D's highlight is set.
$fullRows = @(55, 56, 59, 60, 61, 74, 78, 83)
foreach ($r in $fullRows) {
    $ws.Cells.Item($r, 3).Value = 1          # C: nifti now collected
    $ws.Cells.Item($r, 3).Interior.Pattern = -4142   # clear the "in-progress" highlight

    $ws.Cells.Item($r, 4).Value = 0          # D: recon
    $ws.Cells.Item($r, 4).Interior.ColorIndex = 6    # mark as the new in-progress cell

    $ws.Cells.Item($r, 5).Value = 0          # E: bem_ico4
    $ws.Cells.Item($r, 6).Value = 0          # F: coreg
    $ws.Cells.Item($r, 7).Value = 0          # G: src_surf_ico4
    $ws.Cells.Item($r, 8).Value = 0          # H: src_vol_6.2
}

# Rows 57 and 75 already had C filled in; now D gets filled in too and the
# highlight moves from D to E.
$advanceRows = @(57, 75)
foreach ($r in $advanceRows) {
    $ws.Cells.Item($r, 4).Value = 1          # D: recon
    $ws.Cells.Item($r, 4).Interior.Pattern = -4142   # clear highlight on D

    $ws.Cells.Item($r, 5).Value = 0          # E: bem_ico4 (value unchanged, still 0)
    $ws.Cells.Item($r, 5).Interior.ColorIndex = 6    # move highlight to E
}

# Row 80: AB (rs_preproc) gets newly highlighted, value stays 0.
$ws.Cells.Item(80, 28).Interior.ColorIndex = 6
